$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "struggle"
$ws.Range("C22").Value = -3.367526054382324
$ws.Range("D22").Value = 0.915987193584442
$ws.Range("E22").Value = -1.626443386077881
$ws.Range("F22").Value = -1.0144944190979
$ws.Range("G22").Value = -1.167210817337036
$ws.Range("H22").Value = 0.6551529765129089

$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "struggle"
$ws.Range("C23").Value = -1.118759155273438
$ws.Range("D23").Value = 2.792432069778442
$ws.Range("E23").Value = -4.963344097137451
$ws.Range("F23").Value = -0.3874412775039673
$ws.Range("G23").Value = -0.7050912380218506
$ws.Range("H23").Value = 0.0914770737290382

$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "struggle"
$ws.Range("C24").Value = -1.557756900787354
$ws.Range("D24").Value = -0.5582034587860107
$ws.Range("E24").Value = -0.2619988918304443
$ws.Range("F24").Value = -0.1458440721035003
$ws.Range("G24").Value = -0.3762930035591125
$ws.Range("H24").Value = -0.0704022198915481

$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "struggle"
$ws.Range("C25").Value = 1.715949058532715
$ws.Range("D25").Value = -1.576748490333557
$ws.Range("E25").Value = 5.096891403198242
$ws.Range("F25").Value = 0.2157881408929824
$ws.Range("G25").Value = 0.3240640163421631
$ws.Range("H25").Value = 0.0951422601938247

$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "struggle"
$ws.Range("C26").Value = -0.8243503570556641
$ws.Range("D26").Value = 0.5943599939346313
$ws.Range("E26").Value = 1.927432060241699
$ws.Range("F26").Value = 0.0403171069920063
$ws.Range("G26").Value = 0.1484402567148208
$ws.Range("H26").Value = -0.0852157026529312

$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "struggle"
$ws.Range("C27").Value = -0.2981023788452148
$ws.Range("D27").Value = 1.024843096733093
$ws.Range("E27").Value = 0.8517363667488098
$ws.Range("F27").Value = 0.1014036312699318
$ws.Range("G27").Value = 0.3179553747177124
$ws.Range("H27").Value = 0.0390953756868839

$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "struggle"
$ws.Range("C28").Value = 0.4371089935302734
$ws.Range("D28").Value = 0.3337190449237823
$ws.Range("E28").Value = -0.154114544391632
$ws.Range("F28").Value = 0.052381694316864
$ws.Range("G28").Value = 0.1099557429552078
$ws.Range("H28").Value = 0.0681114718317985

$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "struggle"
$ws.Range("C29").Value = -0.5641984939575195
$ws.Range("D29").Value = -0.3292053341865539
$ws.Range("E29").Value = -0.326197862625122
$ws.Range("F29").Value = 0.0522289797663688
$ws.Range("G29").Value = -0.4196644127368927
$ws.Range("H29").Value = 0.2273945808410644

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "struggle"
$ws.Range("C30").Value = 0.131052017211914
$ws.Range("D30").Value = 0.5107872486114502
$ws.Range("E30").Value = 0.0540084838867187
$ws.Range("F30").Value = 0.0937678143382072
$ws.Range("G30").Value = -0.1565342247486114
$ws.Range("H30").Value = 0.0675006061792373

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "struggle"
$ws.Range("C31").Value = -0.1777238845825195
$ws.Range("D31").Value = 0.4102384448051452
$ws.Range("E31").Value = 0.1352127194404602
$ws.Range("F31").Value = -0.0591012127697467
$ws.Range("G31").Value = 0.0331394411623477
$ws.Range("H31").Value = 0.0291688162833452
